$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) values must stay as exact text (avoid Excel numeric auto-conversion)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.814.99'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.800.30'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.61'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.04'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.798.26'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.33'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.07'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.437.65'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.825.58'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.60'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.783.61'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.10'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '461.46'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.91'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000153'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.50'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.05'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.01'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.00'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.948.42'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.79'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.23'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.66'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.08'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.38'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.996'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '48.08'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '44.05'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.298'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '150.76'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '390.79'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '26.67'

# Columns B, C, E (text / percentage strings)
$ws.Range("E2").Value = '  -1.28%  '
$ws.Range("E3").Value = '  -0.37%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("E5").Value = '  -0.07%  '
$ws.Range("E6").Value = '  -1.50%  '
$ws.Range("E7").Value = '  -0.31%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  -0.34%  '
$ws.Range("E10").Value = '  -0.59%  '
$ws.Range("E11").Value = '  -1.38%  '
$ws.Range("E12").Value = '  -1.04%  '
$ws.Range("E13").Value = '  -3.26%  '
$ws.Range("E14").Value = '  -1.11%  '
$ws.Range("E15").Value = '  -0.69%  '
$ws.Range("E16").Value = '  -0.32%  '
$ws.Range("E17").Value = '  +3.55%  '
$ws.Range("E18").Value = '  -1.31%  '
$ws.Range("E19").Value = '  +1.15%  '
$ws.Range("E20").Value = '  +0.06%  '
$ws.Range("E21").Value = '  -1.21%  '
$ws.Range("E22").Value = '  -8.26%  '
$ws.Range("E23").Value = '  -0.28%  '
$ws.Range("E24").Value = '  -1.78%  '
$ws.Range("E25").Value = '  -0.95%  '
$ws.Range("E26").Value = '  +0.78%  '
$ws.Range("E27").Value = '  -3.37%  '
$ws.Range("E28").Value = '  +0.66%  '
$ws.Range("E29").Value = '  -1.41%  '
$ws.Range("E30").Value = '  -0.54%  '
$ws.Range("E31").Value = '  +0.12%  '
$ws.Range("E32").Value = '  +3.15%  '
$ws.Range("E33").Value = '  -1.39%  '
$ws.Range("E34").Value = '  -1.73%  '
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("E36").Value = '  -1.19%  '
$ws.Range("E37").Value = '  -0.95%  '
$ws.Range("E38").Value = '  -4.13%  '
$ws.Range("E39").Value = '  -0.21%  '
$ws.Range("E40").Value = '  -0.12%  '
$ws.Range("E41").Value = '  -0.52%  '
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("E43").Value = '  -0.01%  '
$ws.Range("E44").Value = '  +2.51%  '
$ws.Range("E45").Value = '  +0.39%  '
$ws.Range("E46").Value = '  -1.82%  '
$ws.Range("E47").Value = '  +2.49%  '
$ws.Range("E48").Value = '  -1.37%  '
$ws.Range("B49").Value = 'Bittensor'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("E49").Value = '  -0.19%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("E50").Value = '  +4.06%  '
$ws.Range("E51").Value = '  -5.71%  '
